$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$rhff  = $wb.Worksheets.Item("RHFF")

# ---------------------------------------------------------------------------
# 1. Add the two new shared strings / fuel rows ("green hydrogen" and
#    "low carbon hydrogen") as new columns (L, M) and new rows (12, 13) on
#    the RHFF matrix, mirroring the existing "To type / From type" table.
# ---------------------------------------------------------------------------

# Header row: new column headers L1/M1
$rhff.Range("L1").Value = "green hydrogen"
$rhff.Range("M1").Value = "low carbon hydrogen"

# New data rows: row 12 = "green hydrogen", row 13 = "low carbon hydrogen"
$rhff.Range("A12").Value = "green hydrogen"
$rhff.Range("A13").Value = "low carbon hydrogen"

# Fill the full matrix values for existing rows 2-11 in the two new columns
# (same pattern as the rest of the table: 1 on the diagonal/"electricity"
# row, 0 everywhere else), and zero-fill the two brand new rows 12-13
# across every column B..M.
$colLetters = @("B","C","D","E","F","G","H","I","J","K","L","M")

for ($r = 2; $r -le 11; $r++) {
    foreach ($col in @("L","M")) {
        $addr = "$col$r"
        if ($r -eq 2) {
            $rhff.Range($addr).Value = 1
        } else {
            $rhff.Range($addr).Value = 0
        }
    }
}

foreach ($r in @(12, 13)) {
    foreach ($col in $colLetters) {
        $rhff.Range("$col$r").Value = 0
    }
}

# Match the right-aligned number style (style index 4) used by the rest of
# the data cells for the newly introduced L:M columns (2-13) and the new
# A12:A13 row labels keep the default style already applied by COM.
$rhff.Range("L2:M13").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 2. Column widths: columns K (11) through the new M (13) end up the same
#    (slightly updated) width, while A..J keep their original width.
# ---------------------------------------------------------------------------
$rhff.Columns.Item(11).ColumnWidth = 14.33
$rhff.Columns.Item(12).ColumnWidth = 14.33
$rhff.Columns.Item(13).ColumnWidth = 14.33

# ---------------------------------------------------------------------------
# 3. Sheet view / selection changes: RHFF becomes the active/selected tab
#    (instead of About), with a selection over the new M2:M13 column. The
#    About sheet keeps its existing B14 selection untouched.
# ---------------------------------------------------------------------------
$rhff.Activate()
$rhff.Range("M2:M13").Select()
